$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated line-flow results (case with 380 kV) for columns C:G across rows 2-25
$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 0.04303700082435569
$block1[0,1] = 0.00657243056490131
$block1[0,2] = 0.07625242059650184
$block1[0,3] = 4.526727007203732
$block1[0,4] = 0.002660522252693241
$block1[1,0] = 0.04263092651696354
$block1[1,1] = 0.00642035348718295
$block1[1,2] = 0.07640180804637442
$block1[1,3] = 4.529732580597994
$block1[1,4] = 0.002665878844378701
$block1[2,0] = 0.04239254619906063
$block1[2,1] = 0.006326336082277173
$block1[2,2] = 0.07651133967862922
$block1[2,3] = 4.533729371605602
$block1[2,4] = 0.002669341828307644
$block1[3,0] = 0.04229817552991477
$block1[3,1] = 0.006287858807517921
$block1[3,2] = 0.07656046474325784
$block1[3,3] = 4.535898411650905
$block1[3,4] = 0.0026707969269183
$block1[4,0] = 0.04228267324837631
$block1[4,1] = 0.006281459668129763
$block1[4,2] = 0.07656889347309637
$block1[4,3] = 4.536291193818755
$block1[4,4] = 0.002671041200859259
$block1[5,0] = 0.04239126223870215
$block1[5,1] = 0.00632581783295727
$block1[5,2] = 0.07651198399999082
$block1[5,3] = 4.533756437221257
$block1[5,4] = 0.002669361274238667
$block1[6,0] = 0.04289472379623405
$block1[6,1] = 0.006520125533471699
$block1[6,2] = 0.07630024080163977
$block1[6,3] = 4.527316402587871
$block1[6,4] = 0.002662333177758626
$block1[7,0] = 0.04396821620876779
$block1[7,1] = 0.00689624182719939
$block1[7,2] = 0.07602579280679933
$block1[7,3] = 4.531792753397227
$block1[7,4] = 0.002649925099024808
$block1[8,0] = 0.04480866988191679
$block1[8,1] = 0.007169850486642204
$block1[8,2] = 0.07590933219838547
$block1[8,3] = 4.545565869748714
$block1[8,4] = 0.002641637078006582
$block1[9,0] = 0.04520209646938866
$block1[9,1] = 0.007293796900704308
$block1[9,2] = 0.07587471431308224
$block1[9,3] = 4.554121212791983
$block1[9,4] = 0.002638044464068759
$block1[10,0] = 0.04535265804949518
$block1[10,1] = 0.007340662713883006
$block1[10,2] = 0.07586423440251444
$block1[10,3] = 4.557691152312742
$block1[10,4] = 0.002636709426928313
$block1[11,0] = 0.04532016193458333
$block1[11,1] = 0.007330572357801657
$block1[11,2] = 0.07586637468289581
$block1[11,3] = 4.556907600361399
$block1[11,4] = 0.00263699582337825
$block1[12,0] = 0.0452144517007298
$block1[12,1] = 0.007297653954884353
$block1[12,2] = 0.07587379949597128
$block1[12,3] = 4.554408290400659
$block1[12,4] = 0.00263793412140754
$block1[13,0] = 0.04514990633324345
$block1[13,1] = 0.007277481510060468
$block1[13,2] = 0.07587868947179466
$block1[13,3] = 4.552920422956049
$block1[13,4] = 0.002638512160917099
$block1[14,0] = 0.04478318045597973
$block1[14,1] = 0.007161740238137426
$block1[14,2] = 0.07591196306871595
$block1[14,3] = 4.54505291563899
$block1[14,4] = 0.002641875426494512
$block1[15,0] = 0.04456103701695469
$block1[15,1] = 0.007090606907183883
$block1[15,2] = 0.07593706999351824
$block1[15,3] = 4.540813606653444
$block1[15,4] = 0.002643984082598525
$block1[16,0] = 0.04443431175102575
$block1[16,1] = 0.007049643696541708
$block1[16,2] = 0.0759532396775171
$block1[16,3] = 4.538590745767266
$block1[16,4] = 0.002645213655329335
$block1[17,0] = 0.04439158487870998
$block1[17,1] = 0.007035765668909733
$block1[17,2] = 0.0759590117258675
$block1[17,3] = 4.537875103591858
$block1[17,4] = 0.002645632844835974
$block1[18,0] = 0.04458457645353064
$block1[18,1] = 0.007098184230075688
$block1[18,2] = 0.07593421847062309
$block1[18,3] = 4.541242580607189
$block1[18,4] = 0.00264375788221932
$block1[19,0] = 0.04524545862485496
$block1[19,1] = 0.007307324747864996
$block1[19,2] = 0.07587154738233615
$block1[19,3] = 4.55513342877444
$block1[19,4] = 0.002637657832013824
$block1[20,0] = 0.0456865821359429
$block1[20,1] = 0.007443604518240221
$block1[20,2] = 0.07584590774060196
$block1[20,3] = 4.566137106515839
$block1[20,4] = 0.002633819127673452
$block1[21,0] = 0.04545030996882815
$block1[21,1] = 0.007370904927595845
$block1[21,2] = 0.07585819399293925
$block1[21,3] = 4.560087765707209
$block1[21,4] = 0.002635854417065952
$block1[22,0] = 0.04457393119687225
$block1[22,1] = 0.007094758734416118
$block1[22,2] = 0.07593550223799639
$block1[22,3] = 4.541047973785453
$block1[22,4] = 0.002643860093417577
$block1[23,0] = 0.04366866230719779
$block1[23,1] = 0.006794991725449506
$block1[23,2] = 0.07608503420182267
$block1[23,3] = 4.528744781275591
$block1[23,4] = 0.00265313569491802
$ws.Range("C2:G25").Value = $block1

# Updated line-flow results for columns I:N across rows 2-25
$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 3.30929981430674
$block2[0,1] = 0.1938379813589499
$block2[0,2] = 3.820272871585871
$block2[0,3] = 0.1514848428649334
$block2[0,4] = 0.9034733304046227
$block2[0,5] = 2.767415841708409
$block2[1,0] = 3.31158005640826
$block2[1,1] = 0.1948076583685889
$block2[1,2] = 3.710516053391189
$block2[1,3] = 0.1521213566546837
$block2[1,4] = 0.885862602215731
$block2[1,5] = 2.795691530960433
$block2[2,0] = 3.314564695323241
$block2[2,1] = 0.1954594255830067
$block2[2,2] = 3.645129515267058
$block2[2,3] = 0.1525458374172608
$block2[2,4] = 0.8754917820256125
$block2[2,5] = 2.813898689983297
$block2[3,0] = 3.316178893397009
$block2[3,1] = 0.1957392175896402
$block2[3,2] = 3.618987733156644
$block2[3,3] = 0.1527272975737191
$block2[3,4] = 0.8713768206866632
$block2[3,5] = 2.821530800279199
$block2[4,0] = 3.316470949080824
$block2[4,1] = 0.1957865343910292
$block2[4,2] = 3.614677328075913
$block2[4,3] = 0.1527579416170735
$block2[4,4] = 0.8707002536143236
$block2[4,5] = 2.822810939851291
$block2[5,0] = 3.314584854447716
$block2[5,1] = 0.1954631414761252
$block2[5,2] = 3.644774918735607
$block2[5,3] = 0.1525482502916464
$block2[5,4] = 0.8754358358057246
$block2[5,5] = 2.814000758896665
$block2[6,0] = 3.309756844210469
$block2[6,1] = 0.1941606322652234
$block2[6,2] = 3.782012428595124
$block2[6,3] = 0.1516973371379731
$block2[6,4] = 0.8973093718093494
$block2[6,5] = 2.776989682811795
$block2[7,0] = 3.312889178343767
$block2[7,1] = 0.1920531955043181
$block2[7,2] = 4.067074535722327
$block2[7,3] = 0.1502950240336869
$block2[7,4] = 0.9437146946457062
$block2[7,5] = 2.711126641780147
$block2[8,0] = 3.322915295883377
$block2[8,1] = 0.1907764930118283
$block2[8,2] = 4.286303782241419
$block2[8,3] = 0.1494261162264223
$block2[8,4] = 0.9799579934529277
$block2[8,5] = 2.666835977108285
$block2[9,0] = 3.329163784529072
$block2[9,1] = 0.190254524022432
$block2[9,2] = 4.388182499378104
$block2[9,3] = 0.149065662739357
$block2[9,4] = 0.9969150249518037
$block2[9,5] = 2.647578172179522
$block2[10,0] = 3.331773341583968
$block2[10,1] = 0.1900653135023447
$block2[10,2] = 4.427071532639047
$block2[10,3] = 0.1489341588546083
$block2[10,4] = 1.00340386025195
$block2[10,5] = 2.640413926715429
$block2[11,0] = 3.331200490278633
$block2[11,1] = 0.1901056877928937
$block2[11,2] = 4.4186822966509
$block2[11,3] = 0.148962258777825
$block2[11,4] = 1.002003367588202
$block2[11,5] = 2.641951166333556
$block2[12,0] = 3.329373592271764
$block2[12,1] = 0.1902387883021781
$block2[12,2] = 4.391375711785372
$block2[12,3] = 0.1490547438707353
$block2[12,4] = 0.9974475108667633
$block2[12,5] = 2.646986193254648
$block2[13,0] = 3.328286282738517
$block2[13,1] = 0.1903214161032061
$block2[13,2] = 4.374689999544728
$block2[13,3] = 0.1491120433394393
$block2[13,4] = 0.9946657182037058
$block2[13,5] = 2.6500870044297
$block2[14,0] = 3.322540961713685
$block2[14,1] = 0.1908117875377258
$block2[14,2] = 4.279689084732354
$block2[14,3] = 0.149450372062752
$block2[14,4] = 0.9788592663665838
$block2[14,5] = 2.668112455200475
$block2[15,0] = 3.319449125370085
$block2[15,1] = 0.1911276705638976
$block2[15,2] = 4.221960221751544
$block2[15,3] = 0.1496668327894106
$block2[15,4] = 0.9692828540274263
$block2[15,5] = 2.679398673865784
$block2[16,0] = 3.317829579702845
$block2[16,1] = 0.1913148940534803
$block2[16,2] = 4.188958459998048
$block2[16,3] = 0.1497946132812302
$block2[16,4] = 0.9638189703839544
$block2[16,5] = 2.685973986443891
$block2[17,0] = 3.317308481139008
$block2[17,1] = 0.1913792357883999
$block2[17,2] = 4.177819368089274
$block2[17,3] = 0.1498384410092441
$block2[17,4] = 0.9619765867392047
$block2[17,5] = 2.688214655721779
$block2[18,0] = 3.319761816643279
$block2[18,1] = 0.1910934713961367
$block2[18,2] = 4.228084610980773
$block2[18,3] = 0.1496434510214542
$block2[18,4] = 0.9702977032070237
$block2[18,5] = 2.67818856415504
$block2[19,0] = 3.329903585135085
$block2[19,1] = 0.190199464254075
$block2[19,2] = 4.399387905726428
$block2[19,3] = 0.1490274433909988
$block2[19,4] = 0.9987838422232471
$block2[19,5] = 2.645503800609642
$block2[20,0] = 3.337950826826003
$block2[20,1] = 0.1896644153309808
$block2[20,2] = 4.51315040785812
$block2[20,3] = 0.1486539375063671
$block2[20,4] = 1.017795078445445
$block2[20,5] = 2.624890181043668
$block2[21,0] = 3.333525780389635
$block2[21,1] = 0.189945478719487
$block2[21,2] = 4.452267779643705
$block2[21,3] = 0.1488506276034478
$block2[21,4] = 1.007612372728033
$block2[21,5] = 2.635823552019364
$block2[22,0] = 3.31961995683362
$block2[22,1] = 0.1911089153393668
$block2[22,2] = 4.225315191139714
$block2[22,3] = 0.1496540115317728
$block2[22,4] = 0.9698387599310792
$block2[22,5] = 2.678735384651731
$block2[23,0] = 3.310688440994312
$block2[23,1] = 0.1925755597268513
$block2[23,2] = 3.988243426539157
$block2[23,3] = 0.150645977718094
$block2[23,4] = 0.9307840377236687
$block2[23,5] = 2.728224913747951
$ws.Range("I2:N25").Value = $block2

